$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 462.66666
$ws.Cells.Item(33, 9).Value = 529.9231
$ws.Cells.Item(33, 11).Value = 529.9231
$ws.Cells.Item(33, 13).Value = -300.9231

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2451.0588
$ws.Cells.Item(98, 9).Value = 2458.16
$ws.Cells.Item(98, 11).Value = 2458.16
$ws.Cells.Item(98, 13).Value = -960.1599999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 2451.0588
$ws.Cells.Item(122, 9).Value = 2458.16
$ws.Cells.Item(122, 11).Value = 7374.48
$ws.Cells.Item(122, 13).Value = -4924.48

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 297698.1
$ws.Cells.Item(132, 9).Value = 297698.1
$ws.Cells.Item(132, 11).Value = 893094.2999999999
$ws.Cells.Item(132, 13).Value = -890564.2999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1809.2241
$ws.Cells.Item(137, 9).Value = 1358.9286
$ws.Cells.Item(137, 11).Value = 4076.7858
$ws.Cells.Item(137, 13).Value = -1526.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12083.35
$ws.Cells.Item(32, 9).Value = 5054.961
$ws.Cells.Item(32, 10).Value = 51910.89
$ws.Cells.Item(32, 11).Value = 5054.961
$ws.Cells.Item(32, 12).Value = 51910.89
$ws.Cells.Item(32, 13).Value = -4767.961
$ws.Cells.Item(32, 14).Value = -52484.89

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 10053.637
$ws.Cells.Item(74, 10).Value = 17700
$ws.Cells.Item(74, 12).Value = 17700
$ws.Cells.Item(74, 14).Value = -19448

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 10053.637
$ws.Cells.Item(77, 10).Value = 17700
$ws.Cells.Item(77, 12).Value = 88500
$ws.Cells.Item(77, 14).Value = -97236

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 3912.2942
$ws.Cells.Item(88, 9).Value = 2356.889
$ws.Cells.Item(88, 11).Value = 2356.889
$ws.Cells.Item(88, 13).Value = -1950.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 3912.2942
$ws.Cells.Item(91, 9).Value = 2356.889
$ws.Cells.Item(91, 11).Value = 2356.889
$ws.Cells.Item(91, 13).Value = -952.8890000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 101008.875
$ws.Cells.Item(110, 9).Value = 1053
$ws.Cells.Item(110, 10).Value = 800700
$ws.Cells.Item(110, 11).Value = 1053
$ws.Cells.Item(110, 12).Value = 800700
$ws.Cells.Item(110, 13).Value = 992
$ws.Cells.Item(110, 14).Value = -804790

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(125, 8).Value = 42485
$ws.Cells.Item(125, 10).Value = 42485
$ws.Cells.Item(125, 12).Value = 42485
$ws.Cells.Item(125, 14).Value = -52325

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1501536.4
$ws.Cells.Item(132, 9).Value = 2031078.8
$ws.Cells.Item(132, 10).Value = 1166.6666
$ws.Cells.Item(132, 11).Value = 6093236.4
$ws.Cells.Item(132, 12).Value = 3499.9998
$ws.Cells.Item(132, 13).Value = -6090706.4
$ws.Cells.Item(132, 14).Value = -8559.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1912.2222
$ws.Cells.Item(99, 9).Value = 1702
$ws.Cells.Item(99, 10).Value = 2332.6667
$ws.Cells.Item(99, 11).Value = 1702
$ws.Cells.Item(99, 12).Value = 2332.6667
$ws.Cells.Item(99, 13).Value = -204
$ws.Cells.Item(99, 14).Value = -5328.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 1912.2222
$ws.Cells.Item(126, 9).Value = 1702
$ws.Cells.Item(126, 10).Value = 2332.6667
$ws.Cells.Item(126, 11).Value = 5106
$ws.Cells.Item(126, 12).Value = 6998.000100000001
$ws.Cells.Item(126, 13).Value = -2636
$ws.Cells.Item(126, 14).Value = -11938.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3107.0667
$ws.Cells.Item(132, 9).Value = 2858.3333
$ws.Cells.Item(132, 10).Value = 3687.4443
$ws.Cells.Item(132, 11).Value = 8574.999899999999
$ws.Cells.Item(132, 12).Value = 11062.3329
$ws.Cells.Item(132, 13).Value = -6044.999899999999
$ws.Cells.Item(132, 14).Value = -16122.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 569.619
$ws.Cells.Item(18, 9).Value = 460.58823
$ws.Cells.Item(18, 10).Value = 1033
$ws.Cells.Item(18, 11).Value = 1381.76469
$ws.Cells.Item(18, 12).Value = 3099
$ws.Cells.Item(18, 13).Value = -1212.76469
$ws.Cells.Item(18, 14).Value = -3437

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 986.2261999999999
$ws.Cells.Item(68, 9).Value = 739.675
$ws.Cells.Item(68, 10).Value = 1210.3636
$ws.Cells.Item(68, 11).Value = 2219.025
$ws.Cells.Item(68, 12).Value = 3631.0908
$ws.Cells.Item(68, 13).Value = -1408.025
$ws.Cells.Item(68, 14).Value = -5253.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 986.2261999999999
$ws.Cells.Item(71, 9).Value = 739.675
$ws.Cells.Item(71, 10).Value = 1210.3636
$ws.Cells.Item(71, 11).Value = 6657.075
$ws.Cells.Item(71, 12).Value = 10893.2724
$ws.Cells.Item(71, 13).Value = -2601.075
$ws.Cells.Item(71, 14).Value = -19005.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(120, 8).Value = 3127.5
$ws.Cells.Item(120, 9).Value = 700
$ws.Cells.Item(120, 10).Value = 5555
$ws.Cells.Item(120, 11).Value = 2100
$ws.Cells.Item(120, 12).Value = 16665
$ws.Cells.Item(120, 13).Value = 2738
$ws.Cells.Item(120, 14).Value = -26341

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 917.9167
$ws.Cells.Item(129, 9).Value = 357.46155
$ws.Cells.Item(129, 10).Value = 1580.2727
$ws.Cells.Item(129, 11).Value = 1072.38465
$ws.Cells.Item(129, 12).Value = 4740.8181
$ws.Cells.Item(129, 13).Value = 3927.61535
$ws.Cells.Item(129, 14).Value = -14740.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1726169.8
$ws.Cells.Item(131, 9).Value = 1257.4166
$ws.Cells.Item(131, 10).Value = 2176147
$ws.Cells.Item(131, 11).Value = 3772.2498
$ws.Cells.Item(131, 12).Value = 6528441
$ws.Cells.Item(131, 13).Value = 1267.7502
$ws.Cells.Item(131, 14).Value = -6538521

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1780.4375
$ws.Cells.Item(102, 9).Value = 1764.75
$ws.Cells.Item(102, 11).Value = 1764.75
$ws.Cells.Item(102, 13).Value = -142.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1869
$ws.Cells.Item(113, 9).Value = 1106.8572
$ws.Cells.Item(113, 11).Value = 1106.8572
$ws.Cells.Item(113, 13).Value = 1063.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 6804.0713
$ws.Cells.Item(122, 9).Value = 7477.909
$ws.Cells.Item(122, 10).Value = 4333.3335
$ws.Cells.Item(122, 11).Value = 22433.727
$ws.Cells.Item(122, 12).Value = 13000.0005
$ws.Cells.Item(122, 13).Value = -19983.727
$ws.Cells.Item(122, 14).Value = -17900.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1757.2727
$ws.Cells.Item(7, 9).Value = 1793.381
$ws.Cells.Item(7, 10).Value = 999
$ws.Cells.Item(7, 11).Value = 1793.381
$ws.Cells.Item(7, 12).Value = 999
$ws.Cells.Item(7, 13).Value = -1681.381
$ws.Cells.Item(7, 14).Value = -1223

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1735.174
$ws.Cells.Item(16, 9).Value = 2102.8823
$ws.Cells.Item(16, 10).Value = 693.3333
$ws.Cells.Item(16, 11).Value = 2102.8823
$ws.Cells.Item(16, 12).Value = 693.3333
$ws.Cells.Item(16, 13).Value = -1932.8823
$ws.Cells.Item(16, 14).Value = -1033.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3100
$ws.Cells.Item(40, 9).Value = 2500
$ws.Cells.Item(40, 10).Value = 4000
$ws.Cells.Item(40, 11).Value = 2500
$ws.Cells.Item(40, 12).Value = 4000
$ws.Cells.Item(40, 13).Value = -2364
$ws.Cells.Item(40, 14).Value = -4272

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 1757.2727
$ws.Cells.Item(126, 9).Value = 1793.381
$ws.Cells.Item(126, 10).Value = 999
$ws.Cells.Item(126, 11).Value = 5380.143
$ws.Cells.Item(126, 12).Value = 2997
$ws.Cells.Item(126, 13).Value = -2910.143
$ws.Cells.Item(126, 14).Value = -7937

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 50001000
$ws.Cells.Item(96, 9).Value = 100000000
$ws.Cells.Item(96, 10).Value = 2000
$ws.Cells.Item(96, 11).Value = 100000000
$ws.Cells.Item(96, 12).Value = 2000
$ws.Cells.Item(96, 13).Value = -99998627
$ws.Cells.Item(96, 14).Value = -4746

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 499.35294
$ws.Cells.Item(113, 9).Value = 379.9
$ws.Cells.Item(113, 10).Value = 670
$ws.Cells.Item(113, 11).Value = 1139.7
$ws.Cells.Item(113, 12).Value = 2010
$ws.Cells.Item(113, 13).Value = 1030.3
$ws.Cells.Item(113, 14).Value = -6350

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1986.4
$ws.Cells.Item(126, 9).Value = 2160.923
$ws.Cells.Item(126, 10).Value = 852
$ws.Cells.Item(126, 11).Value = 6482.768999999999
$ws.Cells.Item(126, 12).Value = 2556
$ws.Cells.Item(126, 13).Value = -4012.768999999999
$ws.Cells.Item(126, 14).Value = -7496
